$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1224
$ws.Range("I43").Value = 1198.6666
$ws.Range("K43").Value = 1198.6666
$ws.Range("M43").Value = -1129.6666
$ws.Range("H58").Value = 1747.2307
$ws.Range("I58").Value = 379.33334
$ws.Range("J58").Value = 4825
$ws.Range("K58").Value = 1138.00002
$ws.Range("L58").Value = 14475
$ws.Range("M58").Value = -988.0000199999999
$ws.Range("N58").Value = -14775
$ws.Range("H80").Value = 1137593
$ws.Range("I80").Value = 1894830.9
$ws.Range("J80").Value = 1736.25
$ws.Range("K80").Value = 5684492.699999999
$ws.Range("L80").Value = 5208.75
$ws.Range("M80").Value = -5683494.699999999
$ws.Range("N80").Value = -7204.75
$ws.Range("H83").Value = 1137593
$ws.Range("I83").Value = 1894830.9
$ws.Range("J83").Value = 1736.25
$ws.Range("K83").Value = 17053478.1
$ws.Range("L83").Value = 15626.25
$ws.Range("M83").Value = -17048486.1
$ws.Range("N83").Value = -25610.25
$ws.Range("H92").Value = 68483.875
$ws.Range("I92").Value = 521.0833
$ws.Range("K92").Value = 521.0833
$ws.Range("M92").Value = 726.9167
$ws.Range("H116").Value = 44964252
$ws.Range("J116").Value = 83343350
$ws.Range("L116").Value = 83343350
$ws.Range("N116").Value = -83350234
$ws.Range("H125").Value = 1550
$ws.Range("J125").Value = 1825
$ws.Range("L125").Value = 16425
$ws.Range("N125").Value = -21345
$ws.Range("H131").Value = 13070
$ws.Range("I131").Value = 3925
$ws.Range("K131").Value = 11775
$ws.Range("M131").Value = -6735
$ws.Range("H137").Value = 52636244
$ws.Range("I137").Value = 100002530
$ws.Range("J137").Value = 7038.4443
$ws.Range("K137").Value = 300007590
$ws.Range("L137").Value = 21115.3329
$ws.Range("M137").Value = -300005040
$ws.Range("N137").Value = -26215.3329
$ws.Range("H138").Value = 3405.8462
$ws.Range("J138").Value = 4234.8887
$ws.Range("L138").Value = 12704.6661
$ws.Range("N138").Value = -22984.6661
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 194916.6
$ws.Range("I32").Value = 215233.94
$ws.Range("K32").Value = 215233.94
$ws.Range("M32").Value = -214946.94
$ws.Range("H61").Value = 1798.3
$ws.Range("J61").Value = 2599.5
$ws.Range("L61").Value = 2599.5
$ws.Range("N61").Value = -3023.5
$ws.Range("H102").Value = 3922.0386
$ws.Range("I102").Value = 1471.8422
$ws.Range("K102").Value = 1471.8422
$ws.Range("M102").Value = 150.1578
$ws.Range("H110").Value = 166668180
$ws.Range("I110").Value = 166668180
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 166668180
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -166666135
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 5953.2
$ws.Range("J122").Value = 7599.8
$ws.Range("L122").Value = 22799.4
$ws.Range("N122").Value = -27699.4
$ws.Range("H136").Value = 1798.3
$ws.Range("J136").Value = 2599.5
$ws.Range("L136").Value = 7798.5
$ws.Range("N136").Value = -12898.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 196
$ws.Range("I22").Value = 170
$ws.Range("K22").Value = 170
$ws.Range("M22").Value = 3
$ws.Range("H94").Value = 2782.7104
$ws.Range("I94").Value = 2860.3
$ws.Range("J94").Value = 2491.75
$ws.Range("K94").Value = 2860.3
$ws.Range("L94").Value = 2491.75
$ws.Range("M94").Value = -2409.3
$ws.Range("N94").Value = -3393.75
$ws.Range("H105").Value = 1821.8823
$ws.Range("J105").Value = 2088.1
$ws.Range("L105").Value = 2088.1
$ws.Range("N105").Value = -5582.1
$ws.Range("H134").Value = 2140.3667
$ws.Range("I134").Value = 1604.3077
$ws.Range("K134").Value = 4812.9231
$ws.Range("M134").Value = -2277.9231
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 48900
$ws.Range("J68").Value = 48900
$ws.Range("L68").Value = 48900
$ws.Range("N68").Value = -50398
$ws.Range("H71").Value = 48900
$ws.Range("J71").Value = 48900
$ws.Range("L71").Value = 146700
$ws.Range("N71").Value = -154188
$ws.Range("H94").Value = 4431.3335
$ws.Range("I94").Value = 4500
$ws.Range("J94").Value = 4397
$ws.Range("K94").Value = 4500
$ws.Range("L94").Value = 4397
$ws.Range("M94").Value = -4049
$ws.Range("N94").Value = -5299
$ws.Range("H107").Value = 900.6087
$ws.Range("I107").Value = 951
$ws.Range("K107").Value = 951
$ws.Range("M107").Value = 969
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 4514.778
$ws.Range("I8").Value = 4514.778
$ws.Range("K8").Value = 13544.334
$ws.Range("M8").Value = -13405.334
$ws.Range("H50").Value = 1187.75
$ws.Range("I50").Value = 626.25
$ws.Range("J50").Value = 1749.25
$ws.Range("K50").Value = 1878.75
$ws.Range("L50").Value = 5247.75
$ws.Range("M50").Value = -1397.75
$ws.Range("N50").Value = -6209.75
$ws.Range("H53").Value = 1187.75
$ws.Range("I53").Value = 626.25
$ws.Range("J53").Value = 1749.25
$ws.Range("K53").Value = 1878.75
$ws.Range("L53").Value = 5247.75
$ws.Range("M53").Value = -1397.75
$ws.Range("N53").Value = -6209.75
$ws.Range("H97").Value = 1095.6154
$ws.Range("I97").Value = 1482.1666
$ws.Range("J97").Value = 764.2857
$ws.Range("K97").Value = 4446.4998
$ws.Range("L97").Value = 2292.8571
$ws.Range("M97").Value = -3950.4998
$ws.Range("N97").Value = -3284.8571
$ws.Range("H121").Value = 93029.46000000001
$ws.Range("I121").Value = 20293.8
$ws.Range("J121").Value = 138489.25
$ws.Range("K121").Value = 60881.39999999999
$ws.Range("L121").Value = 415467.75
$ws.Range("M121").Value = -59571.39999999999
$ws.Range("N121").Value = -418087.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H120").Value = 48976
$ws.Range("J120").Value = 48976
$ws.Range("L120").Value = 48976
$ws.Range("N120").Value = -58652
$ws.Range("H122").Value = 2650.2083
$ws.Range("I122").Value = 1848.2307
$ws.Range("J122").Value = 3598
$ws.Range("K122").Value = 5544.6921
$ws.Range("L122").Value = 10794
$ws.Range("M122").Value = -3094.6921
$ws.Range("N122").Value = -15694
$ws.Range("H126").Value = 6092.4614
$ws.Range("I126").Value = 7565.174
$ws.Range("J126").Value = 3975.4375
$ws.Range("K126").Value = 22695.522
$ws.Range("L126").Value = 11926.3125
$ws.Range("M126").Value = -20225.522
$ws.Range("N126").Value = -16866.3125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 910.26666
$ws.Range("I22").Value = 836.25
$ws.Range("J22").Value = 937.1818
$ws.Range("K22").Value = 836.25
$ws.Range("L22").Value = 937.1818
$ws.Range("M22").Value = -541.25
$ws.Range("N22").Value = -1527.1818
$ws.Range("H27").Value = 910.26666
$ws.Range("I27").Value = 836.25
$ws.Range("J27").Value = 937.1818
$ws.Range("K27").Value = 836.25
$ws.Range("L27").Value = 937.1818
$ws.Range("M27").Value = -729.25
$ws.Range("N27").Value = -1151.1818
$ws.Range("H40").Value = 3596.7058
$ws.Range("I40").Value = 2922.182
$ws.Range("K40").Value = 2922.182
$ws.Range("M40").Value = -2786.182
$ws.Range("H136").Value = 7037.875
$ws.Range("J136").Value = 10329.75
$ws.Range("L136").Value = 30989.25
$ws.Range("N136").Value = -36089.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4363.7144
$ws.Range("J96").Value = 5432.6665
$ws.Range("L96").Value = 5432.6665
$ws.Range("N96").Value = -8178.6665
$ws.Range("H113").Value = 2478.7334
$ws.Range("I113").Value = 1443.7778
$ws.Range("K113").Value = 4332.3334
$ws.Range("M113").Value = -2161.3334
$ws.Range("H132").Value = 280484.56
$ws.Range("I132").Value = 359452.22
$ws.Range("K132").Value = 1078356.66
$ws.Range("M132").Value = -1075826.66
